$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 6502405.5
$ws.Range("I33").Value = 13553.444
$ws.Range("K33").Value = 13553.444
$ws.Range("M33").Value = -13324.444
$ws.Range("H100").Value = 11112992
$ws.Range("I100").Value = 23810926
$ws.Range("J100").Value = 2300
$ws.Range("K100").Value = 23810926
$ws.Range("L100").Value = 2300
$ws.Range("M100").Value = -23810385
$ws.Range("N100").Value = -3382
$ws.Range("H125").Value = 3939.9
$ws.Range("I125").Value = 6216
$ws.Range("J125").Value = 2422.5
$ws.Range("K125").Value = 55944
$ws.Range("L125").Value = 21802.5
$ws.Range("M125").Value = -53484
$ws.Range("N125").Value = -26722.5
$ws.Range("H132").Value = 3003.476
$ws.Range("I132").Value = 2168.2144
$ws.Range("J132").Value = 4674
$ws.Range("K132").Value = 6504.6432
$ws.Range("L132").Value = 14022
$ws.Range("M132").Value = -3974.6432
$ws.Range("N132").Value = -19082
$ws.Range("H138").Value = 1744.3038
$ws.Range("I138").Value = 836.2258
$ws.Range("J138").Value = 2330.7708
$ws.Range("K138").Value = 2508.6774
$ws.Range("L138").Value = 6992.312399999999
$ws.Range("M138").Value = 2631.3226
$ws.Range("N138").Value = -17272.3124

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5953.793
$ws.Range("I32").Value = 4390
$ws.Range("K32").Value = 4390
$ws.Range("M32").Value = -4103
$ws.Range("H61").Value = 5199.1626
$ws.Range("I61").Value = 5857.615
$ws.Range("J61").Value = 4192.1177
$ws.Range("K61").Value = 5857.615
$ws.Range("L61").Value = 4192.1177
$ws.Range("M61").Value = -5645.615
$ws.Range("N61").Value = -4616.1177
$ws.Range("H74").Value = 2119.12
$ws.Range("I74").Value = 1838.5
$ws.Range("J74").Value = 2378.1538
$ws.Range("K74").Value = 1838.5
$ws.Range("L74").Value = 2378.1538
$ws.Range("M74").Value = -964.5
$ws.Range("N74").Value = -4126.1538
$ws.Range("H77").Value = 2119.12
$ws.Range("I77").Value = 1838.5
$ws.Range("J77").Value = 2378.1538
$ws.Range("K77").Value = 9192.5
$ws.Range("L77").Value = 11890.769
$ws.Range("M77").Value = -4824.5
$ws.Range("N77").Value = -20626.769
$ws.Range("H97").Value = 1202.25
$ws.Range("I97").Value = 1048.579
$ws.Range("K97").Value = 1048.579
$ws.Range("M97").Value = -552.579
$ws.Range("H109").Value = 39366.285
$ws.Range("J109").Value = 39366.285
$ws.Range("L109").Value = 39366.285
$ws.Range("N109").Value = -42140.285
$ws.Range("H110").Value = 1378.069
$ws.Range("I110").Value = 1224.2222
$ws.Range("J110").Value = 1629.8182
$ws.Range("K110").Value = 1224.2222
$ws.Range("L110").Value = 1629.8182
$ws.Range("M110").Value = 820.7778000000001
$ws.Range("N110").Value = -5719.8182
$ws.Range("H136").Value = 5199.1626
$ws.Range("I136").Value = 5857.615
$ws.Range("J136").Value = 4192.1177
$ws.Range("K136").Value = 17572.845
$ws.Range("L136").Value = 12576.3531
$ws.Range("M136").Value = -15022.845
$ws.Range("N136").Value = -17676.3531

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 251.2
$ws.Range("I80").Value = 70.5
$ws.Range("J80").Value = 316.9091
$ws.Range("K80").Value = 70.5
$ws.Range("L80").Value = 316.9091
$ws.Range("M80").Value = 927.5
$ws.Range("N80").Value = -2312.9091
$ws.Range("H83").Value = 251.2
$ws.Range("I83").Value = 70.5
$ws.Range("J83").Value = 316.9091
$ws.Range("K83").Value = 352.5
$ws.Range("L83").Value = 1584.5455
$ws.Range("M83").Value = 4639.5
$ws.Range("N83").Value = -11568.5455
$ws.Range("H107").Value = 1258.25
$ws.Range("I107").Value = 1264
$ws.Range("J107").Value = 1255.1538
$ws.Range("K107").Value = 1264
$ws.Range("L107").Value = 1255.1538
$ws.Range("M107").Value = 656
$ws.Range("N107").Value = -5095.1538

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6329.952
$ws.Range("J31").Value = 9995.909
$ws.Range("L31").Value = 9995.909
$ws.Range("N31").Value = -10585.909
$ws.Range("H34").Value = 6329.952
$ws.Range("J34").Value = 9995.909
$ws.Range("L34").Value = 9995.909
$ws.Range("N34").Value = -10399.909
$ws.Range("H58").Value = 1761.5312
$ws.Range("I58").Value = 1487.7
$ws.Range("J58").Value = 2217.9167
$ws.Range("K58").Value = 1487.7
$ws.Range("L58").Value = 2217.9167
$ws.Range("M58").Value = -1284.7
$ws.Range("N58").Value = -2623.9167
$ws.Range("H134").Value = 2596.8113
$ws.Range("I134").Value = 2803.4
$ws.Range("J134").Value = 1961.1538
$ws.Range("K134").Value = 8410.200000000001
$ws.Range("L134").Value = 5883.4614
$ws.Range("M134").Value = -5875.200000000001
$ws.Range("N134").Value = -10953.4614
$ws.Range("H136").Value = 1761.5312
$ws.Range("I136").Value = 1487.7
$ws.Range("J136").Value = 2217.9167
$ws.Range("K136").Value = 4463.1
$ws.Range("L136").Value = 6653.750100000001
$ws.Range("M136").Value = -1913.1
$ws.Range("N136").Value = -11753.7501

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1111677.2
$ws.Range("I113").Value = 1064404.8
$ws.Range("J113").Value = 1429078.1
$ws.Range("K113").Value = 3193214.4
$ws.Range("L113").Value = 4287234.300000001
$ws.Range("M113").Value = -3191044.4
$ws.Range("N113").Value = -4291574.300000001
$ws.Range("H122").Value = 11810.223
$ws.Range("I122").Value = 612.5714
$ws.Range("J122").Value = 51002
$ws.Range("K122").Value = 5513.1426
$ws.Range("L122").Value = 459018
$ws.Range("M122").Value = -3063.1426
$ws.Range("N122").Value = -463918
$ws.Range("H132").Value = 3333.3333
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1612
$ws.Range("I97").Value = 1646.6666
$ws.Range("J97").Value = 1300
$ws.Range("K97").Value = 1646.6666
$ws.Range("L97").Value = 1300
$ws.Range("M97").Value = -1150.6666
$ws.Range("N97").Value = -2292
$ws.Range("H122").Value = 28024024
$ws.Range("I122").Value = 48402796
$ws.Range("K122").Value = 145208388
$ws.Range("M122").Value = -145205938
$ws.Range("H123").Value = 19342.129
$ws.Range("J123").Value = 19686.867
$ws.Range("L123").Value = 19686.867
$ws.Range("N123").Value = -24586.867

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4796609
$ws.Range("I122").Value = 5109455
$ws.Range("J122").Value = 3336661.8
$ws.Range("K122").Value = 15328365
$ws.Range("L122").Value = 10009985.4
$ws.Range("M122").Value = -15325915
$ws.Range("N122").Value = -10014885.4

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 946.64514
$ws.Range("I96").Value = 1072.5
$ws.Range("J96").Value = 717.8182
$ws.Range("K96").Value = 1072.5
$ws.Range("L96").Value = 717.8182
$ws.Range("M96").Value = 300.5
$ws.Range("N96").Value = -3463.8182
$ws.Range("H122").Value = 1188.9231
$ws.Range("I122").Value = 862.625
$ws.Range("J122").Value = 1711
$ws.Range("K122").Value = 2587.875
$ws.Range("L122").Value = 5133
$ws.Range("M122").Value = -137.875
$ws.Range("N122").Value = -10033
